# Auto-generated edit script applying the scraped-schedule refresh
# (new scrape at 10:27:08, 167/19/23 data rows per sheet).
$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2,1).Value = 'Última actualización: 10:27:08'
$ws.Cells.Item(3,1).Value = 'Total filas: 167'
$ws.Cells.Item(47,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(49,3).Value = '225_GOMEZ'
$ws.Cells.Item(54,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(55,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(56,3).Value = '16_SANTA ANA'
$ws.Cells.Item(78,1).Value = '07:18:07'
$ws.Cells.Item(78,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(78,4).Value = 65
$ws.Cells.Item(79,1).Value = '08:16:28'
$ws.Cells.Item(79,3).Value = '215B_EL PATO'
$ws.Cells.Item(79,4).Value = 7
$ws.Cells.Item(118,3).Value = '16_SANTA ANA'
$ws.Cells.Item(119,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(135,1).Value = '10:27:08'
$ws.Cells.Item(135,2).Value = '10:27'
$ws.Cells.Item(135,3).Value = '215A_EL PATO'
$ws.Cells.Item(135,4).Value = 0
$ws.Cells.Item(136,2).Value = '10:32'
$ws.Cells.Item(136,3).Value = '10_OLMOS'
$ws.Cells.Item(136,4).Value = 65
$ws.Cells.Item(137,1).Value = '10:27:08'
$ws.Cells.Item(137,2).Value = '10:33'
$ws.Cells.Item(137,3).Value = '10_OLMOS'
$ws.Cells.Item(137,4).Value = 6
$ws.Cells.Item(138,1).Value = '10:27:08'
$ws.Cells.Item(138,2).Value = '10:35'
$ws.Cells.Item(138,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(138,4).Value = 8
$ws.Cells.Item(139,2).Value = '10:41'
$ws.Cells.Item(139,3).Value = '17_ROMERO'
$ws.Cells.Item(139,4).Value = 74
$ws.Cells.Item(140,1).Value = '10:27:08'
$ws.Cells.Item(140,2).Value = '10:42'
$ws.Cells.Item(140,3).Value = '17_ROMERO'
$ws.Cells.Item(140,4).Value = 15
$ws.Cells.Item(141,2).Value = '10:43'
$ws.Cells.Item(141,3).Value = '14_ABASTO'
$ws.Cells.Item(141,4).Value = 76
$ws.Cells.Item(142,1).Value = '10:27:08'
$ws.Cells.Item(142,2).Value = '10:44'
$ws.Cells.Item(142,3).Value = '14_ABASTO'
$ws.Cells.Item(142,4).Value = 17
$ws.Cells.Item(143,1).Value = '10:27:08'
$ws.Cells.Item(143,2).Value = '10:47'
$ws.Cells.Item(143,3).Value = '16_SANTA ANA'
$ws.Cells.Item(143,4).Value = 20
$ws.Cells.Item(144,1).Value = '10:27:08'
$ws.Cells.Item(144,2).Value = '10:52'
$ws.Cells.Item(144,3).Value = '15_ABASTO'
$ws.Cells.Item(144,4).Value = 25
$ws.Cells.Item(144,5).Value = 'LP1912'
$ws.Cells.Item(145,1).Value = '10:27:08'
$ws.Cells.Item(145,2).Value = '10:53'
$ws.Cells.Item(145,3).Value = '10_OLMOS'
$ws.Cells.Item(145,4).Value = 26
$ws.Cells.Item(145,5).Value = 'LP1912'
$ws.Cells.Item(146,1).Value = '09:27:56'
$ws.Cells.Item(146,2).Value = '10:56'
$ws.Cells.Item(146,3).Value = '27_EL RETIRO'
$ws.Cells.Item(146,4).Value = 89
$ws.Cells.Item(146,5).Value = 'LP1912'
$ws.Cells.Item(147,1).Value = '10:27:08'
$ws.Cells.Item(147,2).Value = '10:57'
$ws.Cells.Item(147,3).Value = '27_EL RETIRO'
$ws.Cells.Item(147,4).Value = 30
$ws.Cells.Item(147,5).Value = 'LP1912'
$ws.Cells.Item(148,1).Value = '09:27:56'
$ws.Cells.Item(148,2).Value = '11:01'
$ws.Cells.Item(148,3).Value = '215C_EL PATO'
$ws.Cells.Item(148,4).Value = 94
$ws.Cells.Item(148,5).Value = 'LP1912'
$ws.Cells.Item(149,1).Value = '10:27:08'
$ws.Cells.Item(149,2).Value = '11:02'
$ws.Cells.Item(149,3).Value = '215C_EL PATO'
$ws.Cells.Item(149,4).Value = 35
$ws.Cells.Item(149,5).Value = 'LP1912'
$ws.Cells.Item(150,1).Value = '10:27:08'
$ws.Cells.Item(150,2).Value = '11:04'
$ws.Cells.Item(150,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(150,4).Value = 37
$ws.Cells.Item(150,5).Value = 'LP1912'
$ws.Cells.Item(151,1).Value = '10:27:08'
$ws.Cells.Item(151,2).Value = '11:05'
$ws.Cells.Item(151,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(151,4).Value = 38
$ws.Cells.Item(151,5).Value = 'LP1912'
$ws.Cells.Item(152,1).Value = '09:27:56'
$ws.Cells.Item(152,2).Value = '11:06'
$ws.Cells.Item(152,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(152,4).Value = 99
$ws.Cells.Item(152,5).Value = 'LP1912'
$ws.Cells.Item(153,1).Value = '10:27:08'
$ws.Cells.Item(153,2).Value = '11:07'
$ws.Cells.Item(153,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(153,4).Value = 40
$ws.Cells.Item(153,5).Value = 'LP1912'
$ws.Cells.Item(154,1).Value = '10:27:08'
$ws.Cells.Item(154,2).Value = '11:12'
$ws.Cells.Item(154,3).Value = '15_ABASTO'
$ws.Cells.Item(154,4).Value = 45
$ws.Cells.Item(154,5).Value = 'LP1912'
$ws.Cells.Item(155,1).Value = '10:27:08'
$ws.Cells.Item(155,2).Value = '11:17'
$ws.Cells.Item(155,3).Value = '16_SANTA ANA'
$ws.Cells.Item(155,4).Value = 50
$ws.Cells.Item(155,5).Value = 'LP1912'
$ws.Cells.Item(156,1).Value = '09:27:56'
$ws.Cells.Item(156,2).Value = '11:19'
$ws.Cells.Item(156,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(156,4).Value = 112
$ws.Cells.Item(156,5).Value = 'LP1912'
$ws.Cells.Item(157,1).Value = '10:27:08'
$ws.Cells.Item(157,2).Value = '11:20'
$ws.Cells.Item(157,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(157,4).Value = 53
$ws.Cells.Item(157,5).Value = 'LP1912'
$ws.Cells.Item(158,1).Value = '09:27:56'
$ws.Cells.Item(158,2).Value = '11:20'
$ws.Cells.Item(158,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(158,4).Value = 113
$ws.Cells.Item(158,5).Value = 'LP1912'
$ws.Cells.Item(159,1).Value = '10:27:08'
$ws.Cells.Item(159,2).Value = '11:21'
$ws.Cells.Item(159,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(159,4).Value = 54
$ws.Cells.Item(159,5).Value = 'LP1912'
$ws.Cells.Item(160,1).Value = '10:27:08'
$ws.Cells.Item(160,2).Value = '11:27'
$ws.Cells.Item(160,3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(160,4).Value = 60
$ws.Cells.Item(160,5).Value = 'LP1912'
$ws.Cells.Item(161,1).Value = '10:27:08'
$ws.Cells.Item(161,2).Value = '11:32'
$ws.Cells.Item(161,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(161,4).Value = 65
$ws.Cells.Item(161,5).Value = 'LP1912'
$ws.Cells.Item(162,1).Value = '10:27:08'
$ws.Cells.Item(162,2).Value = '11:36'
$ws.Cells.Item(162,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(162,4).Value = 69
$ws.Cells.Item(162,5).Value = 'LP1912'
$ws.Cells.Item(163,1).Value = '10:27:08'
$ws.Cells.Item(163,2).Value = '11:42'
$ws.Cells.Item(163,3).Value = '17_ROMERO'
$ws.Cells.Item(163,4).Value = 75
$ws.Cells.Item(163,5).Value = 'LP1912'
$ws.Cells.Item(164,1).Value = '10:27:08'
$ws.Cells.Item(164,2).Value = '11:47'
$ws.Cells.Item(164,3).Value = '10_OLMOS'
$ws.Cells.Item(164,4).Value = 80
$ws.Cells.Item(164,5).Value = 'LP1912'
$ws.Cells.Item(165,1).Value = '10:27:08'
$ws.Cells.Item(165,2).Value = '11:51'
$ws.Cells.Item(165,3).Value = '215B_EL PATO'
$ws.Cells.Item(165,4).Value = 84
$ws.Cells.Item(165,5).Value = 'LP1912'
$ws.Cells.Item(166,1).Value = '10:27:08'
$ws.Cells.Item(166,2).Value = '11:59'
$ws.Cells.Item(166,3).Value = '225_GOMEZ'
$ws.Cells.Item(166,4).Value = 92
$ws.Cells.Item(166,5).Value = 'LP1912'
$ws.Cells.Item(167,1).Value = '10:27:08'
$ws.Cells.Item(167,2).Value = '12:02'
$ws.Cells.Item(167,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(167,4).Value = 95
$ws.Cells.Item(167,5).Value = 'LP1912'
$ws.Cells.Item(168,1).Value = '10:27:08'
$ws.Cells.Item(168,2).Value = '12:07'
$ws.Cells.Item(168,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(168,4).Value = 100
$ws.Cells.Item(168,5).Value = 'LP1912'
$ws.Cells.Item(169,1).Value = '10:27:08'
$ws.Cells.Item(169,2).Value = '12:14'
$ws.Cells.Item(169,3).Value = '17_ROMERO'
$ws.Cells.Item(169,4).Value = 107
$ws.Cells.Item(169,5).Value = 'LP1912'
$ws.Cells.Item(170,1).Value = '10:27:08'
$ws.Cells.Item(170,2).Value = '12:19'
$ws.Cells.Item(170,3).Value = '14_ABASTO'
$ws.Cells.Item(170,4).Value = 112
$ws.Cells.Item(170,5).Value = 'LP1912'
$ws.Cells.Item(171,1).Value = '10:27:08'
$ws.Cells.Item(171,2).Value = '12:21'
$ws.Cells.Item(171,3).Value = '215A_EL PATO'
$ws.Cells.Item(171,4).Value = 114
$ws.Cells.Item(171,5).Value = 'LP1912'
$ws.Cells.Item(172,1).Value = '10:27:08'
$ws.Cells.Item(172,2).Value = '12:21'
$ws.Cells.Item(172,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(172,4).Value = 114
$ws.Cells.Item(172,5).Value = 'LP1912'

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2,1).Value = 'Última actualización: 10:27:08'
$ws.Cells.Item(3,1).Value = 'Total filas: 19'
$ws.Cells.Item(20,1).Value = '10:27:08'
$ws.Cells.Item(20,2).Value = '10:27'
$ws.Cells.Item(20,3).Value = '215A_EL PATO'
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(21,1).Value = '09:27:56'
$ws.Cells.Item(21,2).Value = '11:01'
$ws.Cells.Item(21,3).Value = '215C_EL PATO'
$ws.Cells.Item(21,4).Value = 94
$ws.Cells.Item(21,5).Value = 'LP1912'
$ws.Cells.Item(22,1).Value = '10:27:08'
$ws.Cells.Item(22,2).Value = '11:02'
$ws.Cells.Item(22,3).Value = '215C_EL PATO'
$ws.Cells.Item(22,4).Value = 35
$ws.Cells.Item(22,5).Value = 'LP1912'
$ws.Cells.Item(23,1).Value = '10:27:08'
$ws.Cells.Item(23,2).Value = '11:51'
$ws.Cells.Item(23,3).Value = '215B_EL PATO'
$ws.Cells.Item(23,4).Value = 84
$ws.Cells.Item(23,5).Value = 'LP1912'
$ws.Cells.Item(24,1).Value = '10:27:08'
$ws.Cells.Item(24,2).Value = '12:21'
$ws.Cells.Item(24,3).Value = '215A_EL PATO'
$ws.Cells.Item(24,4).Value = 114
$ws.Cells.Item(24,5).Value = 'LP1912'

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2,1).Value = 'Última actualización: 10:27:08'
$ws.Cells.Item(3,1).Value = 'Total filas: 23'
$ws.Cells.Item(25,1).Value = '10:27:08'
$ws.Cells.Item(25,2).Value = '10:54'
$ws.Cells.Item(25,3).Value = '215A_LA PLATA'
$ws.Cells.Item(25,4).Value = 27
$ws.Cells.Item(25,5).Value = 'L6173'
$ws.Cells.Item(26,1).Value = '09:27:56'
$ws.Cells.Item(26,2).Value = '11:13'
$ws.Cells.Item(26,3).Value = '215C_LA PLATA'
$ws.Cells.Item(26,4).Value = 106
$ws.Cells.Item(26,5).Value = 'L6203'
$ws.Cells.Item(27,1).Value = '10:27:08'
$ws.Cells.Item(27,2).Value = '11:14'
$ws.Cells.Item(27,3).Value = '215C_LA PLATA'
$ws.Cells.Item(27,4).Value = 47
$ws.Cells.Item(27,5).Value = 'L6203'
$ws.Cells.Item(28,1).Value = '10:27:08'
$ws.Cells.Item(28,2).Value = '12:04'
$ws.Cells.Item(28,3).Value = '215A_LA PLATA'
$ws.Cells.Item(28,4).Value = 97
$ws.Cells.Item(28,5).Value = 'L6173'
